$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows
$ws.Range("F2").Value = 8
$ws.Range("F5").Value = -3
$ws.Range("F7").Value = -11
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = 8
$ws.Range("F12").Value = -4
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = 2
